# "Handling of additional FieldAttrTypes"
# Qvx Reader handles more FieldAttrTypes, FixPointDecimals. Rewriting of
# Reader documentation.
#
# This touches only the "Tasks" sheet (sheet1): the "Reader (Monica)" /
# "Writer (Matthew)" column headers are renamed to plain "Reader" /
# "Writer"; the previously red/not-done rows for UNKNOWN, ASCII, INTEGER,
# REAL, FIX and MONEY are marked green/done; and the DATE, TIME, TIMESTAMP
# and INTERVAL rows (still red/not-done) now note that a "Reg-exp function"
# is needed under the Reader column. Focus also moves from the
# "Matt Planning" sheet back to the "Tasks" sheet.

$wb = $excel.ActiveWorkbook
$wsTasks = $wb.Worksheets.Item("Tasks")

# Rename the header cells (was "Reader (Monica)" / "Writer (Matthew)")
$wsTasks.Range("B1").Value = "Reader"
$wsTasks.Range("C1").Value = "Writer"

# UNKNOWN, ASCII, INTEGER, REAL, FIX, MONEY (rows 2-7): Reader column goes
# from red ("not done") to green ("done")
$wsTasks.Range("B2:B7").Interior.Color = 5287936

# DATE, TIME, TIMESTAMP, INTERVAL (rows 8-11): still red/not-done, but now
# annotated with what's needed
$wsTasks.Range("B8:B11").Value = "Reg-exp function"

# Switch the active sheet/selection back to "Tasks" (was "Matt Planning")
$wsTasks.Activate()
$wsTasks.Range("B7").Select()
